# Fill in the attendance marks for the 26th (column G) on the November
# attendance sheet, and leave the view scrolled/selected where the user
# was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows ---------------------------------------------------
# Row 10 ("Date"): the 3rd date column (G) is the 26th.
$ws.Range("G10").Value = 26

# Row 11 ("Time"): G11 was an empty cell merged/styled differently from
# the other date columns; match D11/E11/F11's format (wrapped text font)
# before writing the same time-slot text into it, exactly like the other
# populated date columns.
$ws.Range("F11").Copy()
$ws.Range("G11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G11").Value = "09:00 AM To 12:00 PM"

# --- "Total Classes" row -------------------------------------------
$ws.Range("G12").Value = 3

# --- Per-student attendance marks for the 26th (column G) ----------
$ws.Range("G14").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("G16").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("G20").Value = 3
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 3
$ws.Range("G23").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("G25").Value = 0

# --- View state: scroll so row 11 is at the top, selection on G24 --
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G24").Select()
